# 项目计划表.xlsx — add the "第七周周四" (2018.10.11) weekly block and mark
# the previous week's checklist as completed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Column C ("完成情况") for the 2018.10.10 block is now "已完成" ---
foreach ($r in 23..28) {
    $ws.Cells.Item($r, 3).Value = "已完成"
}

# --- 2. The summary line for that same block now carries the full note ---
$ws.Range("A29").Value = "总结：任务基本已完成，细节有待改善。"

# --- 3. Append a brand-new weekly block in rows 31-40 -------------------

# Row 31: date banner (will be merged A31:D31)
$ws.Range("A31").Value = "日期：2018.10.11 第七周周四"

# Row 32: column headers
$ws.Range("A32").Value = "组员"
$ws.Range("B32").Value = "计划内容"
$ws.Range("C32").Value = "完成情况"
$ws.Range("D32").Value = "备注"

# Rows 33-38: the same six team members, plan/status left blank for now
$ws.Range("A33").Value = "练富珊"
$ws.Range("A34").Value = "黄成志"
$ws.Range("A35").Value = "黄皓燊"
$ws.Range("A36").Value = "郑嘉蔚"
$ws.Range("A37").Value = "郑瑞贤"
$ws.Range("A38").Value = "辛伟达"

# Row 39: summary placeholder (will be merged A39:D40), row 40 stays blank
$ws.Range("A39").Value = "总结："

# --- 4. Merge the banner/summary rows, matching the earlier blocks ------
$ws.Range("A31:D31").Merge()
$ws.Range("A39:D40").Merge()

# --- 5. Re-apply the same visual styling used by the existing blocks ----

# Date banner row -> same look as A1/A11/A21 (bold 10pt, centered, boxed)
$banner = $ws.Range("A31:D31")
$banner.Borders.LineStyle = 1
$banner.Borders.Weight = 2
$banner.Font.Bold = $true
$banner.Font.Size = 10
$banner.HorizontalAlignment = -4108
$banner.VerticalAlignment = -4108

# Column-header row -> same look as A2:D2 / A12:D12 / A22:D22
$header = $ws.Range("A32:D32")
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2
$header.Font.Bold = $true
$header.Font.Size = 10
$header.VerticalAlignment = -4108

# Team-member data rows -> same look as A3:D8 / A13:D18 / A23:D28
$data = $ws.Range("A33:D38")
$data.Borders.LineStyle = 1
$data.Borders.Weight = 2
$data.Font.Bold = $false
$data.Font.Size = 11
$data.VerticalAlignment = -4108

# Summary rows -> same look as A9:D10 / A19:D20 / A29:D30
$summary = $ws.Range("A39:D40")
$summary.Borders.LineStyle = 1
$summary.Borders.Weight = 2
$summary.Font.Bold = $false
$summary.Font.Size = 11
$summary.HorizontalAlignment = -4131
$summary.VerticalAlignment = -4108

# --- 6. Leave the selection where the author left it when saving --------
$ws.Range("D34").Select()
